$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused rows (23-43) so the sheet shrinks back to 22 rows.
$ws.Range("A23:B43").EntireRow.Delete()

# Normalize every data row (2-22) to the same staff name / branch values.
$ws.Range("A2:A22").Value = "Alex"
$ws.Range("B2:B22").Value = "Utrecht"

# Match the saved selection/active-cell state from the diff.
$ws.Range("B12").Select()
